$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pipes")

# Copy header style (bold, border, alignment) from F1 onto the new header cells
$ws.Range("F1").Copy()
$ws.Range("G1:M1").PasteSpecial(-4122)

$headers = @("Nvol", "Direction", "uMax", "MMgas", "Tgas", "xCH4", "Z")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 7 + $i
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 7).Value = 2
}
